$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header / title area
# ------------------------------------------------------------------
$ws.Range("G1").Value = "Freeze plug"
$ws.Range("M1").Value = "High res"

# New "Freeze plug" parameter table header row (row 2)
$ws.Range("H2").Value = "Degrees"
$ws.Range("I2").Value = "Wall width"
$ws.Range("J2").Value = "DT"
$ws.Range("K2").Value = "Freezing"
$ws.Range("L2").Value = "Melting"
$ws.Range("M2").Value = "Freezing"
$ws.Range("N2").Value = "Melting"
$ws.Range("H2:N2").Font.Bold = $true

# ------------------------------------------------------------------
# New "Freeze plug" parameter table data (rows 3-15)
# ------------------------------------------------------------------
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 40

$ws.Range("J5").Value = 80

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 40

$ws.Range("J7").Value = 80

$ws.Range("H8").Value = 60
$ws.Range("I8").Value = 2.5
$ws.Range("J8").Value = 40

$ws.Range("J9").Value = 80

$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 40

$ws.Range("J11").Value = 80

$ws.Range("H12").Value = 90
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = 40

$ws.Range("J13").Value = 80

$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 40

$ws.Range("J15").Value = 80

# Result placeholder cells - mark the ones that still need a value with
# a red "to fill in" highlight; leave the not-yet-applicable ones unfilled.
$redCells = @("L4","L5","N5","L6","K7","L7","L8","K9","L9","L10","K11","L11","L12","K13","L13","L14","K15","L15")
foreach ($addr in $redCells) {
    $ws.Range($addr).Interior.Color = 192
}

$noFillCells = @("H3","I3","K3","L3","K4","K6","K8","H14","K10","K12","H15","I15","K14")
foreach ($addr in $noFillCells) {
    $ws.Range($addr).Interior.ColorIndex = -4142
}

$ws.Range("K5").Interior.Color = $ws.Range("B6").Interior.Color
$ws.Range("M5").Interior.Color = $ws.Range("B6").Interior.Color

# ------------------------------------------------------------------
# New row of plot links (insert before the existing "Plots" row so the
# rest of the plot-link table shifts down by one row)
# ------------------------------------------------------------------
$ws.Rows("23:23").Insert()
$ws.Range("A24").ClearContents()

$ws.Range("A23").Value = "Plots"
$ws.Range("C23").Value = "Octadecane"
$ws.Range("G23").Value = "Natural convection"
$ws.Range("J23").Value = "Stefan"
$ws.Range("L23").Value = "Freeze plug"
$ws.Range("A23").Font.Bold = $true
$ws.Range("C23").Font.Bold = $true
$ws.Range("G23").Font.Bold = $true
$ws.Range("J23").Font.Bold = $true
$ws.Range("L23").Font.Bold = $true

$ws.Range("L24").Value = "Benchmark with Shafer"
$ws.Range("N24").Value = "10 cm, 5 K"

# ------------------------------------------------------------------
# Misc view / layout tweaks
# ------------------------------------------------------------------
$ws.Columns("G:G").ColumnWidth = 12.33203125
$ws.Range("M11").Select()
